$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 1: subject count headers for columns B:E
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 (CON): updated meanEMG values
$ws.Range("B2").Value = 169.76763831018951
$ws.Range("C2").Value = 136.6188627315039
$ws.Range("D2").Value = 171.01743350577439
$ws.Range("E2").Value = 134.77546343603575

# Row 3 (STR): updated meanEMG values
$ws.Range("B3").Value = 153.98837884296645
$ws.Range("C3").Value = 128.60867508149147
$ws.Range("D3").Value = 158.61882371519948
$ws.Range("E3").Value = 132.45973953741387

# Update the active selection to match the narrower range used after edit
$ws.Range("B1:E3").Select()
